# Updated results to April 22, removed France, changed file names of output

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-7 (row 1 is the header, unchanged).
# Columns: A=Country, B=Region, C=Date(serial), D=CFR2, E=DiffIT, F=AgeCompIT,
#          G=RateCompIT, H=relAgeIT, I=relRateIT

$data = @(
    @("SouthKorea", "All", 43943, 0.02226482139517487, 0.1077415102689016, 0.07287847691095452, 0.03486303335794706, 0.676419670831272, 0.323580329168728),
    @("China",      "All", 43872, 0.02290025071633238, 0.1071060809477441, 0.05419755578227433, 0.05290852516546973, 0.5060175416997729, 0.4939824583002271),
    @("Germany",    "All", 43943, 0.03807511017180802, 0.09193122149226841, 0.04858740299581205, 0.04334381849645638, 0.5285190624808388, 0.4714809375191612),
    @("USA",        "All", 43943, 0.05598346774242222, 0.07402286392165422, 0.07038013310537267, 0.003642730816281557, 0.9507891126701472, 0.04921088732985287),
    @("Spain",      "All", 43943, 0.1040122728013148, 0.02599405886276164, 0.01328374519896857, 0.01271031366379308, 0.5110300499472394, 0.4889699500527606),
    @("Italy",      "All", 43943, 0.1300063316640764, 0, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]

    if ($vals.Count -gt 7) {
        $ws.Cells.Item($row, 8).Value = $vals[7]
        $ws.Cells.Item($row, 9).Value = $vals[8]
    } else {
        $ws.Cells.Item($row, 8).ClearContents()
        $ws.Cells.Item($row, 9).ClearContents()
    }
}

# Rows 8 and 9 (formerly Spain/Italy "All" rows) are no longer needed since
# the table now ends at row 7 - delete them (shifting cells up) so the
# used range shrinks back to A1:I7.
$ws.Range("A8:I9").EntireRow.Delete()
